$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price/Volume columns so numeric-looking strings
# (e.g. "0.00001096", "307.45") are preserved verbatim as text, matching
# the original inlineStr cell content instead of being coerced to numbers.
$textCells = @("D2","D3","D5","D6","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D24","D27","D28","D29","D30","D31","D32","D33","D34","D35","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "24.013.78"
$ws.Range("E2").Value = "  +15.83%  "
$ws.Range("D3").Value = "1.665.57"
$ws.Range("E3").Value = "  +10.32%  "
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").Value = "307.45"
$ws.Range("E5").Value = "  +10.07%  "
$ws.Range("D6").Value = "0.9964"
$ws.Range("E6").Value = "  +3.57%  "
$ws.Range("D7").Value = "0.3720"
$ws.Range("E7").Value = "  +4.46%  "
$ws.Range("D8").Value = "0.3449"
$ws.Range("E8").Value = "  +11.02%  "
$ws.Range("D9").Value = "47.76"
$ws.Range("E9").Value = "  +20.71%  "
$ws.Range("D10").Value = "1.167"
$ws.Range("E10").Value = "  +6.42%  "
$ws.Range("D11").Value = "0.07244"
$ws.Range("E11").Value = "  +7.93%  "
$ws.Range("D12").Value = "0.9983"
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("D13").Value = "20.43"
$ws.Range("E13").Value = "  +11.17%  "
$ws.Range("D14").Value = "6.017"
$ws.Range("E14").Value = "  +8.64%  "
$ws.Range("D15").Value = "6.722"
$ws.Range("E15").Value = "  +7.92%  "
$ws.Range("D16").Value = "1.659.44"
$ws.Range("E16").Value = "  +11.17%  "
$ws.Range("D17").Value = "0.00001096"
$ws.Range("E17").Value = "  +7.03%  "
$ws.Range("D18").Value = "0.9960"
$ws.Range("E18").Value = "  +3.60%  "
$ws.Range("D19").Value = "0.06712"
$ws.Range("E19").Value = "  +11.57%  "
$ws.Range("D20").Value = "82.07"
$ws.Range("E20").Value = "  +17.64%  "
$ws.Range("D21").Value = "16.39"
$ws.Range("E21").Value = "  +11.84%  "
$ws.Range("E22").Value = "  +10.10%  "
$ws.Range("E23").Value = "  +7.06%  "
$ws.Range("D24").Value = "23.975.96"
$ws.Range("E24").Value = "  +15.92%  "
$ws.Range("E25").Value = "  +3.63%  "
$ws.Range("E26").Value = "  -7.85%  "
$ws.Range("D27").Value = "2.663"
$ws.Range("E27").Value = "  +25.27%  "
$ws.Range("D28").Value = "152.21"
$ws.Range("E28").Value = "  +4.07%  "
$ws.Range("D29").Value = "19.49"
$ws.Range("E29").Value = "  +11.78%  "
$ws.Range("D30").Value = "1.841.87"
$ws.Range("E30").Value = "  +11.70%  "
$ws.Range("D31").Value = "126.99"
$ws.Range("E31").Value = "  +9.84%  "
$ws.Range("D32").Value = "6.368"
$ws.Range("E32").Value = "  +26.72%  "
$ws.Range("D33").Value = "4.068"
$ws.Range("E33").Value = "  +2.24%  "
$ws.Range("D34").Value = "0.9775"
$ws.Range("E34").Value = "  +19.54%  "
$ws.Range("D35").Value = "1.730"
$ws.Range("E35").Value = "  +20.44%  "
$ws.Range("E36").Value = "  +5.45%  "
$ws.Range("D37").Value = "12.29"
$ws.Range("E37").Value = "  +17.74%  "
$ws.Range("D38").Value = "8.962"
$ws.Range("E38").Value = "  +17.78%  "
$ws.Range("D39").Value = "5.314"
$ws.Range("E39").Value = "  +10.76%  "
$ws.Range("D40").Value = "0.06345"
$ws.Range("E40").Value = "  +10.27%  "
$ws.Range("D41").Value = "1.289"
$ws.Range("E41").Value = "  +7.11%  "
$ws.Range("D42").Value = "0.02313"
$ws.Range("E42").Value = "  +13.35%  "
$ws.Range("D43").Value = "0.2076"
$ws.Range("E43").Value = "  +11.25%  "
$ws.Range("D44").Value = "0.6096"
$ws.Range("E44").Value = "  +15.41%  "
$ws.Range("D45").Value = "0.9955"
$ws.Range("E45").Value = "  +3.43%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "13.29"
$ws.Range("E46").Value = "  +9.65%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "3.808"
$ws.Range("E47").Value = "  +7.83%  "
$ws.Range("D48").Value = "0.5940"
$ws.Range("E48").Value = "  +13.45%  "
$ws.Range("D49").Value = "127.11"
$ws.Range("E49").Value = "  +5.96%  "
$ws.Range("D50").Value = "2.002"
$ws.Range("E50").Value = "  +8.95%  "
$ws.Range("E51").Value = "  +9.58%  "
